# Append the 2025-11-16 profit row (row 91) to the bottom of the data,
# following the same layout as the existing rows: column A holds the
# date as literal text "MM/DD/YYYY", column B holds the numeric profit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A must stay literal text (not an auto-converted date serial),
# matching every other date cell above it, so force a text format on
# the new cell before typing the date-like string into it.
$ws.Range("A91").NumberFormat = "@"
$ws.Range("A91").Value = "11/16/2025"

$ws.Range("B91").Value = 8762.84
